$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.193.89'
$ws.Range("E2").Value = '  +1.77%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.097.72'
$ws.Range("E3").Value = '  +0.69%  '

# Row 4
$ws.Range("E4").Value = '  -0.10%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '539.07'
$ws.Range("E5").Value = '  -1.51%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.19'
$ws.Range("E6").Value = '  -0.20%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.07%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.093.46'
$ws.Range("E8").Value = '  +0.77%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.496'
$ws.Range("E9").Value = '  +1.02%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.55'
$ws.Range("E10").Value = '  +0.57%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.157'
$ws.Range("E11").Value = '  -0.44%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.459'
$ws.Range("E12").Value = '  +0.05%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000228'
$ws.Range("E13").Value = '  +5.18%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.87'
$ws.Range("E14").Value = '  +0.36%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.592.29'
$ws.Range("E15").Value = '  +0.57%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.155.07'
$ws.Range("E16").Value = '  +1.58%  '

# Row 17
$ws.Range("E17").Value = '  +0.96%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.095.09'
$ws.Range("E18").Value = '  +0.49%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.71'
$ws.Range("E19").Value = '  +1.26%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '487.35'
$ws.Range("E20").Value = '  -2.84%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.47'
$ws.Range("E21").Value = '  +0.75%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.705'
$ws.Range("E22").Value = '  +0.81%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.15'
$ws.Range("E23").Value = '  -0.02%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '80.01'
$ws.Range("E24").Value = '  +3.11%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.29'
$ws.Range("E25").Value = '  +0.32%  '

# Row 26
$ws.Range("E26").Value = '  +0.24%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.72'
$ws.Range("E27").Value = '  -0.18%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.23'
$ws.Range("E28").Value = '  -0.80%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.998'
$ws.Range("E29").Value = '  -0.26%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '26.38'
$ws.Range("E30").Value = '  +0.30%  '

# Row 31
$ws.Range("E31").Value = '  +2.84%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.90'
$ws.Range("E32").Value = '  -1.76%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.39'
$ws.Range("E33").Value = '  -4.55%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '56.87'
$ws.Range("E34").Value = '  -1.46%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '505.80'
$ws.Range("E35").Value = '  -2.87%  '

# Row 36
$ws.Range("E36").Value = '  +5.26%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.05'
$ws.Range("E37").Value = '  +2.02%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.249.30'
$ws.Range("E38").Value = '  +6.13%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0400'
$ws.Range("E39").Value = '  -0.76%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0797'
$ws.Range("E40").Value = '  +1.02%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.119'
$ws.Range("E41").Value = '  -0.16%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.17'
$ws.Range("E42").Value = '  +1.01%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.67'
$ws.Range("E43").Value = '  +0.59%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.257'
$ws.Range("E44").Value = '  +0.55%  '

# Row 45
$ws.Range("E45").Value = '  -0.02%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.07'
$ws.Range("E46").Value = '  +1.58%  '

# Row 47
$ws.Range("B47").Value = 'PEPE'
$ws.Range("C47").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₃0537'
$ws.Range("E47").Value = '  +7.11%  '

# Row 48
$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '122.52'
$ws.Range("E48").Value = '  +0.88%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '24.67'
$ws.Range("E49").Value = '  +1.72%  '

# Row 50
$ws.Range("B50").Value = 'CoreDAO'
$ws.Range("C50").Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.52'
$ws.Range("E50").Value = '  -0.75%  '

# Row 51
$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.109'
$ws.Range("E51").Value = '  +2.44%  '
